# Gatchina area 14-16 (added)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1) Update existing rows 51-53 ("companies" / N column and R51)
# ---------------------------------------------------------------
$ws.Range("N51").Value = 11911
$ws.Range("N52").Value = 11312
$ws.Range("N53").Value = 10838
$ws.Range("R51").Value = 1340

# Give N51:N53 the highlighted ("companies corrected") look already
# used elsewhere in the sheet (yellow fill + centered) by copying the
# format of an existing highlighted cell (P6) onto them.
$hl = $ws.Range("P6")
$hl.Copy()
$ws.Range("N51:N53").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------
# 2) Add the three new "Гатчинский МР" rows: 54 (2016), 55 (2015)
#    and 56 (2014)
# ---------------------------------------------------------------

# Row 54 - clone formatting from row 53 first, then overwrite values
$ws.Range("A53:U53").Copy()
$ws.Range("A54:U54").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A54").Value = "Гатчинский МР"
$ws.Range("B54").Value = 2016
$ws.Range("C54").Value = 245.60599999999999
$ws.Range("D54").Value = 36.497999999999998
$ws.Range("E54").Value = 1156
$ws.Range("F54").Value = 37629.5
$ws.Range("G54").Value = "???"
$ws.Range("H54").Value = "???"
$ws.Range("I54").Value = "???"
$ws.Range("J54").Value = "???"
$ws.Range("K54").Value = "???"
$ws.Range("L54").Value = 7320.8
$ws.Range("M54").Value = "???"
$ws.Range("N54").Value = 10551
$ws.Range("O54").Value = 60117.7
$ws.Range("P54").Value = 1270
$ws.Range("Q54").Value = 138.69999999999999
$ws.Range("R54").Value = 1684
$ws.Range("S54").Formula = "= 13562423.7 / 1000"
$ws.Range("T54").Formula = "= 128209.1 / 1000"
$ws.Range("U54").Value = 767

# Row 55
$ws.Range("A53:U53").Copy()
$ws.Range("A55:U55").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A55").Value = "Гатчинский МР"
$ws.Range("B55").Value = 2015
$ws.Range("C55").Value = 246
$ws.Range("D55").Value = 36.811
$ws.Range("E55").Value = 1380
$ws.Range("F55").Value = 34278.9
$ws.Range("G55").Value = "???"
$ws.Range("H55").Value = "???"
$ws.Range("I55").Value = "???"
$ws.Range("J55").Value = "???"
$ws.Range("K55").Value = "???"
$ws.Range("L55").Value = 7689.2
$ws.Range("M55").Value = "???"
$ws.Range("N55").Value = 10455
$ws.Range("O55").Value = 57475.199999999997
$ws.Range("P55").Value = 2335.9
$ws.Range("Q55").Value = 181.1
$ws.Range("R55").Value = 2530
$ws.Range("S55").Value = 31739
$ws.Range("U55").Value = 942

# T55 uses a distinct "no fill / centered" style that does not exist
# yet in the workbook - create it by toggling the interior pattern
# (forces a brand new xf record) after centering the cell.
$ws.Range("T55").Value = 722
$ws.Range("T55").HorizontalAlignment = -4108
$ws.Range("T55").Interior.Pattern = 1
$ws.Range("T55").Interior.Pattern = -4142

# Give N54/N55 the same highlighted look as N51:N53
$hl.Copy()
$ws.Range("N54:N55").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("N54").Value = 10551
$ws.Range("N55").Value = 10455

# Row 56 - note columns H:K stay completely empty on this row, and L56
# keeps the workbook's default (unstyled) formatting, so copy the
# formats in two pieces (A:G and L:U) leaving H:K untouched.
$ws.Range("A53:G53").Copy()
$ws.Range("A56:G56").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("L53:U53").Copy()
$ws.Range("L56:U56").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("L56").ClearFormats()

$ws.Range("A56").Value = "Гатчинский МР"
$ws.Range("B56").Value = 2014
$ws.Range("C56").Value = 246.2
$ws.Range("D56").Value = 37.902999999999999
$ws.Range("E56").Value = 763
$ws.Range("F56").Value = 32674.7
$ws.Range("G56").Value = 25.9
$ws.Range("L56").Formula = "= 4558593 / 1000"
$ws.Range("M56").Value = "???"
$ws.Range("N56").Value = "???"
$ws.Range("O56").Formula = "= 41975492 / 1000"
$ws.Range("P56").Formula = "= 3280900 / 1000"
$ws.Range("Q56").Value = 123.8
$ws.Range("R56").Value = "???"
$ws.Range("S56").Formula = "= 10145823.2 / 1000"
$ws.Range("T56").Formula = "= 119837 / 1000"
$ws.Range("U56").Value = 3215

# ---------------------------------------------------------------
# 3) Update the view state (scrolled down a bit further, selection
#    moved to F57 like in the source workbook)
# ---------------------------------------------------------------
$ws.Range("F57").Select()
